# Generate Report for handback
#
# For both language sheets ("zh-cn" and "de-de") this:
#   1. Flips the Status column (B) from "Ready for handoff" to
#      "Handed back: in sync with en-US" for the two real source rows.
#   2. Fills in the "Latest Target File" (E) and "Latest Handback File" (F)
#      columns for those two rows, re-using the same hyperlink targets as
#      the "Source File Name" (A) and "Latest Handoff File" (C) columns.
#   3. Records the actual handback timestamp in the
#      "Latest Handback DateTime" column (G), replacing the
#      "0001-01-01 00:00:00" placeholder.

function Get-HyperlinkAddress($ws, $addr) {
    $found = ""
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $found = $h.Address
        }
    }
    return $found
}

function Update-HandbackSheet($wb, $sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Same display text reused by both rows' new E/F cells.
    $sourceMdDisplay = $ws.Range("A2").Text
    $xlfDisplay = $ws.Range("C2").Text

    # Hyperlink targets to reuse.
    $mdUrl = Get-HyperlinkAddress $ws '$A$2'
    $xlfUrl = Get-HyperlinkAddress $ws '$C$2'

    # 1. Status column -> handed back.
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # 2. Latest Target File / Latest Handback File columns.
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, "", "", $sourceMdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl, "", "", $xlfDisplay)
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, "", "", $sourceMdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl, "", "", $xlfDisplay)

    # 3. Latest Handback DateTime column -> actual timestamp.
    $ws.Range("G2").Value = $handbackDateTime
    $ws.Range("G3").Value = $handbackDateTime
}

$wb = $excel.ActiveWorkbook

Update-HandbackSheet $wb "zh-cn" "2016-01-26 06:59:15"
Update-HandbackSheet $wb "de-de" "2016-01-26 06:59:36"
